$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# Fill in previously-empty row 21 with DORN_nohints data
$ws.Range("A21").Value = "DORN_nohints"
$ws.Range("B21").Value = 0.97945953669757002
$ws.Range("C21").Value = 0.994678041453117
$ws.Range("D21").Value = 0.99806731906079904
$ws.Range("E21").Value = 0.118581589964839
$ws.Range("F21").Value = 0.30914099379332699
$ws.Range("G21").Value = 0.084022855021325601
$ws.Range("H21").Value = 0.031475905720026003
$ws.Range("I21").Value = 0.038626020417079
$ws.Range("J21").Value = 0.105912321823124

# Update A31 from "DORN" to "DORN_nohints"
$ws.Range("A31").Value = "DORN_nohints"

# Update selection
$ws.Range("A32").Select()
